$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.9155336666666667
$ws.Range("H2").Value = 2.746601
$ws.Range("I2").Value = 0.1890240037548773
$ws.Range("J2").Value = 0.1951995261655112
$ws.Range("M2").Value = 6.038588
$ws.Range("N2").Value = 18.115764
$ws.Range("O2").Value = 0.4832124175152646
$ws.Range("P2").Value = 0.5250116720691949
$ws.Range("Q2").Value = 5.528530613129333
$ws.Range("R2").Value = 49.756775518164
$ws.Range("S2").Value = 0.0913387458228087
$ws.Range("T2").Value = 0.1024820296192696
$ws.Range("G3").Value = 0.9155336666666667
$ws.Range("H3").Value = 2.746601
$ws.Range("I3").Value = 0.1890240037548773
$ws.Range("J3").Value = 0.1951995261655112
$ws.Range("O3").Value = 0.157682833439183
$ws.Range("P3").Value = 0.1713228489992161
$ws.Range("Q3").Value = 1.804081062974667
$ws.Range("R3").Value = 16.236729566772
$ws.Range("S3").Value = 0.02980584050008782
$ws.Range("T3").Value = 0.03344213894597239
$ws.Range("G4").Value = 0.9155336666666667
$ws.Range("H4").Value = 2.746601
$ws.Range("I4").Value = 0.1890240037548773
$ws.Range("J4").Value = 0.1951995261655112
$ws.Range("M4").Value = 0.6106236666666667
$ws.Range("N4").Value = 1.831871
$ws.Range("O4").Value = 0.04886257154189607
$ws.Range("P4").Value = 0.05308932357062435
$ws.Range("Q4").Value = 0.5590465244967778
$ws.Range("R4").Value = 5.031418720471001
$ws.Range("S4").Value = 0.009236198906608322
$ws.Range("T4").Value = 0.01036301080543338
$ws.Range("G5").Value = 0.9155336666666667
$ws.Range("H5").Value = 2.746601
$ws.Range("I5").Value = 0.1890240037548773
$ws.Range("J5").Value = 0.1951995261655112
$ws.Range("M5").Value = 2.98482
$ws.Range("N5").Value = 5.96964
$ws.Range("O5").Value = 0.2388475729836035
$ws.Range("P5").Value = 0.173005713590172
$ws.Range("Q5").Value = 2.73270319894
$ws.Range("R5").Value = 16.39621919364
$ws.Range("S5").Value = 0.04514792453249599
$ws.Range("T5").Value = 0.03377063331672771
$ws.Range("G6").Value = 0.9155336666666667
$ws.Range("H6").Value = 2.746601
$ws.Range("I6").Value = 0.1890240037548773
$ws.Range("J6").Value = 0.1951995261655112
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.892201
$ws.Range("N6").Value = 2.676603
$ws.Range("O6").Value = 0.07139460452005281
$ws.Range("P6").Value = 0.07757044177079274
$ws.Range("Q6").Value = 0.8168400529336667
$ws.Range("R6").Value = 7.351560476403001
$ws.Range("S6").Value = 0.01349529399287644
$ws.Range("T6").Value = 0.01514171347810812
$ws.Range("I7").Value = 0.7156204889943075
$ws.Range("J7").Value = 0.7390002200311341
$ws.Range("M7").Value = 6.038588
$ws.Range("N7").Value = 18.115764
$ws.Range("O7").Value = 0.4832124175152646
$ws.Range("P7").Value = 0.5250116720691949
$ws.Range("Q7").Value = 20.93030357095867
$ws.Range("R7").Value = 188.372732138628
$ws.Range("S7").Value = 0.3457967065103952
$ws.Range("T7").Value = 0.3879837411780486
$ws.Range("I8").Value = 0.7156204889943075
$ws.Range("J8").Value = 0.7390002200311341
$ws.Range("O8").Value = 0.157682833439183
$ws.Range("P8").Value = 0.1713228489992161
$ws.Range("S8").Value = 0.1128410663717561
$ws.Range("T8").Value = 0.1266076231067814
$ws.Range("I9").Value = 0.7156204889943075
$ws.Range("J9").Value = 0.7390002200311341
$ws.Range("M9").Value = 0.6106236666666667
$ws.Range("N9").Value = 1.831871
$ws.Range("O9").Value = 0.04886257154189607
$ws.Range("P9").Value = 0.05308932357062435
$ws.Range("Q9").Value = 2.116478009585222
$ws.Range("R9").Value = 19.048302086267
$ws.Range("S9").Value = 0.034967057340331
$ws.Range("T9").Value = 0.03923302179999547
$ws.Range("I10").Value = 0.7156204889943075
$ws.Range("J10").Value = 0.7390002200311341
$ws.Range("M10").Value = 2.98482
$ws.Range("N10").Value = 5.96964
$ws.Range("O10").Value = 0.2388475729836035
$ws.Range("P10").Value = 0.173005713590172
$ws.Range("Q10").Value = 10.34566171838
$ws.Range("R10").Value = 62.07397031028
$ws.Range("S10").Value = 0.1709242169736299
$ws.Range("T10").Value = 0.1278512604097805
$ws.Range("I11").Value = 0.7156204889943075
$ws.Range("J11").Value = 0.7390002200311341
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.892201
$ws.Range("N11").Value = 2.676603
$ws.Range("O11").Value = 0.07139460452005281
$ws.Range("P11").Value = 0.07757044177079274
$ws.Range("Q11").Value = 3.092451045892334
$ws.Range("R11").Value = 27.832059413031
$ws.Range("S11").Value = 0.05109144179819539
$ws.Range("T11").Value = 0.05732457353652811
$ws.Range("F12").Value = 0.6666666666666666
$ws.Range("G12").Value = 0.002153333333333333
$ws.Range("H12").Value = 0.00646
$ws.Range("I12").Value = 0.0004445840747369229
$ws.Range("J12").Value = 0.0004591088909634862
$ws.Range("M12").Value = 6.038588
$ws.Range("N12").Value = 18.115764
$ws.Range("O12").Value = 0.4832124175152646
$ws.Range("P12").Value = 0.5250116720691949
$ws.Range("Q12").Value = 0.01300309282666667
$ws.Range("R12").Value = 0.11702783544
$ws.Range("S12").Value = 0.0002148285455424156
$ws.Range("T12").Value = 0.0002410375265065736
$ws.Range("F13").Value = 0.6666666666666666
$ws.Range("G13").Value = 0.002153333333333333
$ws.Range("H13").Value = 0.00646
$ws.Range("I13").Value = 0.0004445840747369229
$ws.Range("J13").Value = 0.0004591088909634862
$ws.Range("O13").Value = 0.157682833439183
$ws.Range("P13").Value = 0.1713228489992161
$ws.Range("Q13").Value = 0.004243195013333334
$ws.Range("R13").Value = 0.03818875512
$ws.Range("S13").Value = 0.0000701032766064555
$ws.Range("T13").Value = 0.0000786558432007349
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.002153333333333333
$ws.Range("H14").Value = 0.00646
$ws.Range("I14").Value = 0.0004445840747369229
$ws.Range("J14").Value = 0.0004591088909634862
$ws.Range("M14").Value = 0.6106236666666667
$ws.Range("N14").Value = 1.831871
$ws.Range("O14").Value = 0.04886257154189607
$ws.Range("P14").Value = 0.05308932357062435
$ws.Range("Q14").Value = 0.001314876295555556
$ws.Range("R14").Value = 0.01183388666
$ws.Range("S14").Value = 0.00002172352115822056
$ws.Range("T14").Value = 0.00002437378046651101
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.002153333333333333
$ws.Range("H15").Value = 0.00646
$ws.Range("I15").Value = 0.0004445840747369229
$ws.Range("J15").Value = 0.0004591088909634862
$ws.Range("M15").Value = 2.98482
$ws.Range("N15").Value = 5.96964
$ws.Range("O15").Value = 0.2388475729836035
$ws.Range("P15").Value = 0.173005713590172
$ws.Range("Q15").Value = 0.0064273124
$ws.Range("R15").Value = 0.0385638744
$ws.Range("S15").Value = 0.000106187827238075
$ws.Range("T15").Value = 0.0000794284612967304
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.002153333333333333
$ws.Range("H16").Value = 0.00646
$ws.Range("I16").Value = 0.0004445840747369229
$ws.Range("J16").Value = 0.0004591088909634862
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 0.892201
$ws.Range("N16").Value = 2.676603
$ws.Range("O16").Value = 0.07139460452005281
$ws.Range("P16").Value = 0.07757044177079274
$ws.Range("Q16").Value = 0.001921206153333334
$ws.Range("R16").Value = 0.01729085538
$ws.Range("S16").Value = 0.00003174090419175621
$ws.Range("T16").Value = 0.00003561327949293634
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 0.459699
$ws.Range("H17").Value = 0.9193979999999999
$ws.Range("I17").Value = 0.09491092317607834
$ws.Range("J17").Value = 0.06534114491239122
$ws.Range("M17").Value = 6.038588
$ws.Range("N17").Value = 18.115764
$ws.Range("O17").Value = 0.4832124175152646
$ws.Range("P17").Value = 0.5250116720691949
$ws.Range("Q17").Value = 2.775932865012
$ws.Range("R17").Value = 16.655597190072
$ws.Range("S17").Value = 0.04586213663651837
$ws.Range("T17").Value = 0.03430486374537008
$ws.Range("E18").Value = 2
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = 0.459699
$ws.Range("H18").Value = 0.9193979999999999
$ws.Range("I18").Value = 0.09491092317607834
$ws.Range("J18").Value = 0.06534114491239122
$ws.Range("O18").Value = 0.157682833439183
$ws.Range("P18").Value = 0.1713228489992161
$ws.Range("Q18").Value = 0.9058479122759999
$ws.Range("R18").Value = 5.435087473655999
$ws.Range("S18").Value = 0.01496582329073266
$ws.Range("T18").Value = 0.0111944311032615
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = 0.459699
$ws.Range("H19").Value = 0.9193979999999999
$ws.Range("I19").Value = 0.09491092317607834
$ws.Range("J19").Value = 0.06534114491239122
$ws.Range("M19").Value = 0.6106236666666667
$ws.Range("N19").Value = 1.831871
$ws.Range("O19").Value = 0.04886257154189607
$ws.Range("P19").Value = 0.05308932357062435
$ws.Range("Q19").Value = 0.280703088943
$ws.Range("R19").Value = 1.684218533658
$ws.Range("S19").Value = 0.00463759177379853
$ws.Range("T19").Value = 0.003468917184728993
$ws.Range("E20").Value = 2
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = 0.459699
$ws.Range("H20").Value = 0.9193979999999999
$ws.Range("I20").Value = 0.09491092317607834
$ws.Range("J20").Value = 0.06534114491239122
$ws.Range("M20").Value = 2.98482
$ws.Range("N20").Value = 5.96964
$ws.Range("O20").Value = 0.2388475729836035
$ws.Range("P20").Value = 0.173005713590172
$ws.Range("Q20").Value = 1.37211876918
$ws.Range("R20").Value = 5.488475076719999
$ws.Range("S20").Value = 0.02266924365023956
$ws.Range("T20").Value = 0.01130439140236708
$ws.Range("E21").Value = 2
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = 0.459699
$ws.Range("H21").Value = 0.9193979999999999
$ws.Range("I21").Value = 0.09491092317607834
$ws.Range("J21").Value = 0.06534114491239122
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.892201
$ws.Range("N21").Value = 2.676603
$ws.Range("O21").Value = 0.07139460452005281
$ws.Range("P21").Value = 0.07757044177079274
$ws.Range("Q21").Value = 0.410143907499
$ws.Range("R21").Value = 2.460863444994
$ws.Range("S21").Value = 0.006776127824789227
$ws.Range("T21").Value = 0.005068541476663573
